# Update gh-pages to output generated at 456a3b4
$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("G2").Value = 79
$ws1.Range("F6").Value = 4437
$ws1.Range("F7").Value = 333

# Sheet "演出" (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 16

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("G2").Value = 79
$ws4.Range("F6").Value = 4437
$ws4.Range("F7").Value = 333
$ws4.Range("F11").Value = 16
